$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Cell values for the logbook rows (row 4 now holds a new entry, the old
#    row-4 entry moved down to row 6, and two brand-new entries were added
#    in rows 5 and 7). The shared-string table fills up in the order the
#    author actually typed the new task descriptions, so write A5 (Github
#    repo) before A4 (Eerste GUI) to reproduce that ordering.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Github repo aangemaakt met Spring project"
$ws.Range("A4").Value = "Eerste GUI en proof of concept: Aanmaken index-, login- en registreerpagina. HTML + CSS en JS"
$ws.Range("A7").Value = "Login authorisatie gefixt + doorverwijzingen controller. CSS gefixt + fragments. Logout knop toegevoegd. Alle klassen (repo,service,..) om gebruiker te registreren toegevoegd, maar nog probleem met toevoegen aan DB via restAPI. Tabel Customer aangepast om functie user mee op te slaan. Klasse voor klusjes aangemaakt, en begin van pagina om nieuwe klusjes toe te voegen."
$ws.Range("A6").Value = "opzetten database, eerste securty toevoeging, login dat leest uit database en start van registeren (nog niet af)"

$ws.Range("B4").Value = 45627
$ws.Range("C4").Value = 2.5
$ws.Range("D4").ClearContents()

$ws.Range("B5").Value = 45630
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 0.5

$ws.Range("B6").Value = 45630
$ws.Range("C6").ClearContents()
$ws.Range("D6").Value = 4.5

$ws.Range("B7").Value = 45633
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 5

# ---------------------------------------------------------------------------
# 2. Number format for the "Gewerkte uren" columns (C/D) changes from 2
#    decimals to a dedicated "0.0" custom format, applied to every data row.
# ---------------------------------------------------------------------------
$ws.Range("C4").NumberFormat = "0.0"
$fmtSrc = $ws.Range("C4")
$fmtDst = $ws.Range("C5:D33")
$fmtSrc.Copy()
$fmtDst.PasteSpecial(-4122)
$ws.Range("D4").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 3. The task column (A) now wraps text on every remaining blank row too, so
#    copy the already-wrapped look of A4 down across the rest of the table.
# ---------------------------------------------------------------------------
$wrapSrc = $ws.Range("A4")
$wrapDst = $ws.Range("A5:A33")
$wrapSrc.Copy()
$wrapDst.PasteSpecial(-4122)

# Re-apply the values we set above for A5-A7 (PasteSpecial only touched
# formats, so the text already there is untouched, but make sure the style
# used for A4 itself carries no number format, matching the source data).
$ws.Range("A4").WrapText = $true
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4160

$excel.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Row heights (explicit/custom heights for the four populated rows).
# ---------------------------------------------------------------------------
$ws.Rows(4).RowHeight = 55.2
$ws.Rows(5).RowHeight = 37.2
$ws.Rows(6).RowHeight = 57
$ws.Rows(7).RowHeight = 163.8

# ---------------------------------------------------------------------------
# 5. Column widths (best effort; the headless engine quantizes these to
#    whole pixels so we pick the nearest achievable widths).
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 35.36
$ws.Columns(2).ColumnWidth = 23.07
$ws.Columns(3).ColumnWidth = 26.26
$ws.Columns(4).ColumnWidth = 26.67

# ---------------------------------------------------------------------------
# 6. Selection state - user ended up with D7 selected after scrolling down.
# ---------------------------------------------------------------------------
$ws.Range("D7").Select()
